# Apply the "Added the rest of the team photos" edit:
#  - User Stories sheet: drop the C-scores for rows 6-8 and 11-13 from 5 -> 4,
#    and clear out rows 14-19 (student rows that hadn't actually been filled
#    in yet for that batch).
#  - Project Management sheet: rescale a few rubric rows (4, 8, 9, 10) down,
#    and blank out row 11 (Burndown Chart) entirely.

$wb = $excel.ActiveWorkbook

# ---- "User Stories" sheet -------------------------------------------------
$wsUser = $wb.Worksheets.Item("User Stories")

$wsUser.Range("C6").Value = 4
$wsUser.Range("C7").Value = 4
$wsUser.Range("C8").Value = 4

$wsUser.Range("C11").Value = 4
$wsUser.Range("C12").Value = 4
$wsUser.Range("C13").Value = 4

# Rows 14-19 (student id + grade columns) get cleared back to blank.
$wsUser.Range("A14:C19").ClearContents()

# ---- "Project Management" sheet -------------------------------------------
$wsPM = $wb.Worksheets.Item("Project Management")

$wsPM.Range("C4:G4").Value = 3

$wsPM.Range("C8:G8").Value = 3

$wsPM.Range("C9").Value = 3
$wsPM.Range("D9").Value = 3
$wsPM.Range("E9").Value = 3
$wsPM.Range("F9").Value = 4
$wsPM.Range("G9").Value = 3

$wsPM.Range("C10:G10").Value = 3

$wsPM.Range("C11:G11").ClearContents()
